$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing year column (Q) into the new column (R)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Set the new values
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 42.9

# Update the active selection like the original workbook did (one row below headers)
$ws.Range("R9").Select()
